$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.783.40"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.375.25"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.26"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.03"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.374.10"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.947.18"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.71"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.371.61"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("E17").Value = "  -3.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.937.71"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.80"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.71"
$ws.Range("E20").Value = "  -3.63%  "
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "372.27"
$ws.Range("E22").Value = "  -1.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.510.35"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.546"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.91"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.177"
$ws.Range("E28").Value = "  +10.55%  "
$ws.Range("E29").Value = "  -4.27%  "
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.31"
$ws.Range("E31").Value = "  -3.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.13"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.24"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("E36").Value = "  -4.36%  "
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.77"
$ws.Range("E38").Value = "  -0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.72"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0755"
$ws.Range("E40").Value = "  -3.35%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.772"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.97"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.68"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("E46").Value = "  -5.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.530.08"
$ws.Range("E47").Value = "  +8.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.27"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.75"
$ws.Range("E49").Value = "  -1.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.40"
$ws.Range("E50").Value = "  +4.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0257"
$ws.Range("E51").Value = "  -1.36%  "
